# chore: update Sheets via scheduled runner
# Refreshes cached market-board figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) on a handful of leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW
# sheets, as produced by the scheduled data-refresh job.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(34, 8).Value = 5045.5835
$ws.Cells.Item(34, 9).Value = 2149.8
$ws.Cells.Item(34, 10).Value = 19524.5
$ws.Cells.Item(34, 11).Value = 2149.8
$ws.Cells.Item(34, 12).Value = 19524.5
$ws.Cells.Item(34, 13).Value = -1946.8
$ws.Cells.Item(34, 14).Value = -19930.5
$ws.Cells.Item(36, 8).Value = 5045.5835
$ws.Cells.Item(36, 9).Value = 2149.8
$ws.Cells.Item(36, 10).Value = 19524.5
$ws.Cells.Item(36, 11).Value = 2149.8
$ws.Cells.Item(36, 12).Value = 19524.5
$ws.Cells.Item(36, 13).Value = -1434.8
$ws.Cells.Item(36, 14).Value = -20954.5
$ws.Cells.Item(51, 8).Value = 3661.8333
$ws.Cells.Item(51, 9).Value = 3533.3845
$ws.Cells.Item(51, 10).Value = 3813.6365
$ws.Cells.Item(51, 11).Value = 3533.3845
$ws.Cells.Item(51, 12).Value = 3813.6365
$ws.Cells.Item(51, 13).Value = -3049.3845
$ws.Cells.Item(51, 14).Value = -4781.636500000001
$ws.Cells.Item(113, 8).Value = 4049.125
$ws.Cells.Item(113, 9).Value = 2696
$ws.Cells.Item(113, 10).Value = 4664.1816
$ws.Cells.Item(113, 11).Value = 2696
$ws.Cells.Item(113, 12).Value = 4664.1816
$ws.Cells.Item(113, 13).Value = 558
$ws.Cells.Item(113, 14).Value = -11172.1816
$ws.Cells.Item(129, 8).Value = 682.8125
$ws.Cells.Item(129, 9).Value = 300.8
$ws.Cells.Item(129, 10).Value = 1319.5
$ws.Cells.Item(129, 11).Value = 902.4000000000001
$ws.Cells.Item(129, 12).Value = 3958.5
$ws.Cells.Item(129, 13).Value = 4097.6
$ws.Cells.Item(129, 14).Value = -13958.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 391133.38
$ws.Cells.Item(61, 9).Value = 257326.98
$ws.Cells.Item(61, 10).Value = 914194.75
$ws.Cells.Item(61, 11).Value = 257326.98
$ws.Cells.Item(61, 12).Value = 914194.75
$ws.Cells.Item(61, 13).Value = -257114.98
$ws.Cells.Item(61, 14).Value = -914618.75
$ws.Cells.Item(88, 8).Value = 2428.28
$ws.Cells.Item(88, 9).Value = 2721.4285
$ws.Cells.Item(88, 10).Value = 2055.182
$ws.Cells.Item(88, 11).Value = 2721.4285
$ws.Cells.Item(88, 12).Value = 2055.182
$ws.Cells.Item(88, 13).Value = -2315.4285
$ws.Cells.Item(88, 14).Value = -2867.182
$ws.Cells.Item(91, 8).Value = 2428.28
$ws.Cells.Item(91, 9).Value = 2721.4285
$ws.Cells.Item(91, 10).Value = 2055.182
$ws.Cells.Item(91, 11).Value = 2721.4285
$ws.Cells.Item(91, 12).Value = 2055.182
$ws.Cells.Item(91, 13).Value = -1317.4285
$ws.Cells.Item(91, 14).Value = -4863.182
$ws.Cells.Item(136, 8).Value = 391133.38
$ws.Cells.Item(136, 9).Value = 257326.98
$ws.Cells.Item(136, 10).Value = 914194.75
$ws.Cells.Item(136, 11).Value = 771980.9400000001
$ws.Cells.Item(136, 12).Value = 2742584.25
$ws.Cells.Item(136, 13).Value = -769430.9400000001
$ws.Cells.Item(136, 14).Value = -2747684.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4943.3887
$ws.Cells.Item(86, 9).Value = 6128.5454
$ws.Cells.Item(86, 10).Value = 3081
$ws.Cells.Item(86, 11).Value = 6128.5454
$ws.Cells.Item(86, 12).Value = 3081
$ws.Cells.Item(86, 13).Value = -5005.5454
$ws.Cells.Item(86, 14).Value = -5327
$ws.Cells.Item(89, 8).Value = 4943.3887
$ws.Cells.Item(89, 9).Value = 6128.5454
$ws.Cells.Item(89, 10).Value = 3081
$ws.Cells.Item(89, 11).Value = 30642.727
$ws.Cells.Item(89, 12).Value = 15405
$ws.Cells.Item(89, 13).Value = -25026.727
$ws.Cells.Item(89, 14).Value = -26637

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2608.9436
$ws.Cells.Item(31, 9).Value = 1844
$ws.Cells.Item(31, 10).Value = 4205.3477
$ws.Cells.Item(31, 11).Value = 1844
$ws.Cells.Item(31, 12).Value = 4205.3477
$ws.Cells.Item(31, 13).Value = -1549
$ws.Cells.Item(31, 14).Value = -4795.3477
$ws.Cells.Item(34, 8).Value = 2608.9436
$ws.Cells.Item(34, 9).Value = 1844
$ws.Cells.Item(34, 10).Value = 4205.3477
$ws.Cells.Item(34, 11).Value = 1844
$ws.Cells.Item(34, 12).Value = 4205.3477
$ws.Cells.Item(34, 13).Value = -1642
$ws.Cells.Item(34, 14).Value = -4609.3477

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(35, 8).Value = 750
$ws.Cells.Item(35, 9).Value = 1000
$ws.Cells.Item(35, 10).Value = 500
$ws.Cells.Item(35, 11).Value = 3000
$ws.Cells.Item(35, 12).Value = 1500
$ws.Cells.Item(35, 13).Value = -2712
$ws.Cells.Item(35, 14).Value = -2076
$ws.Cells.Item(57, 8).Value = 0
$ws.Cells.Item(57, 9).Value = 0
$ws.Cells.Item(57, 10).Value = 0
$ws.Cells.Item(57, 11).Value = 0
$ws.Cells.Item(57, 12).Value = 0
$ws.Cells.Item(57, 14).ClearContents()
$ws.Cells.Item(58, 8).Value = 2781.5293
$ws.Cells.Item(58, 9).Value = 1000
$ws.Cells.Item(58, 10).Value = 3019.0667
$ws.Cells.Item(58, 11).Value = 3000
$ws.Cells.Item(58, 12).Value = 9057.2001
$ws.Cells.Item(58, 13).Value = -2872
$ws.Cells.Item(58, 14).Value = -9313.2001
$ws.Cells.Item(74, 8).Value = 2296.4167
$ws.Cells.Item(74, 9).Value = 895.4
$ws.Cells.Item(74, 10).Value = 3297.1428
$ws.Cells.Item(74, 11).Value = 2686.2
$ws.Cells.Item(74, 12).Value = 9891.428400000001
$ws.Cells.Item(74, 13).Value = -1625.2
$ws.Cells.Item(74, 14).Value = -12013.4284
$ws.Cells.Item(77, 8).Value = 2296.4167
$ws.Cells.Item(77, 9).Value = 895.4
$ws.Cells.Item(77, 10).Value = 3297.1428
$ws.Cells.Item(77, 11).Value = 8058.599999999999
$ws.Cells.Item(77, 12).Value = 29674.2852
$ws.Cells.Item(77, 13).Value = -2754.599999999999
$ws.Cells.Item(77, 14).Value = -40282.2852
$ws.Cells.Item(100, 8).Value = 5173.7144
$ws.Cells.Item(100, 9).Value = 4525
$ws.Cells.Item(100, 10).Value = 5281.8335
$ws.Cells.Item(100, 11).Value = 13575
$ws.Cells.Item(100, 12).Value = 15845.5005
$ws.Cells.Item(100, 13).Value = -12764
$ws.Cells.Item(100, 14).Value = -17467.5005
$ws.Cells.Item(124, 8).Value = 1890
$ws.Cells.Item(124, 9).Value = 280
$ws.Cells.Item(124, 10).Value = 3500
$ws.Cells.Item(124, 11).Value = 840
$ws.Cells.Item(124, 12).Value = 10500
$ws.Cells.Item(124, 13).Value = 4070
$ws.Cells.Item(124, 14).Value = -20320

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(31, 8).Value = 732.4
$ws.Cells.Item(31, 9).Value = 732.4
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 732.4
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = -440.4
$ws.Cells.Item(37, 8).Value = 732.4
$ws.Cells.Item(37, 9).Value = 732.4
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 11).Value = 732.4
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 13).Value = -455.4
$ws.Cells.Item(70, 8).Value = 3883.3333
$ws.Cells.Item(70, 9).Value = 3920
$ws.Cells.Item(70, 10).Value = 3778.5715
$ws.Cells.Item(70, 11).Value = 3920
$ws.Cells.Item(70, 12).Value = 3778.5715
$ws.Cells.Item(70, 13).Value = -3650
$ws.Cells.Item(70, 14).Value = -4318.5715
$ws.Cells.Item(73, 8).Value = 3883.3333
$ws.Cells.Item(73, 9).Value = 3920
$ws.Cells.Item(73, 10).Value = 3778.5715
$ws.Cells.Item(73, 11).Value = 3920
$ws.Cells.Item(73, 12).Value = 3778.5715
$ws.Cells.Item(73, 13).Value = -2984
$ws.Cells.Item(73, 14).Value = -5650.5715
$ws.Cells.Item(80, 8).Value = 5224.8213
$ws.Cells.Item(80, 9).Value = 6352.647
$ws.Cells.Item(80, 10).Value = 3481.818
$ws.Cells.Item(80, 11).Value = 6352.647
$ws.Cells.Item(80, 12).Value = 3481.818
$ws.Cells.Item(80, 13).Value = -5354.647
$ws.Cells.Item(80, 14).Value = -5477.818
$ws.Cells.Item(83, 8).Value = 5224.8213
$ws.Cells.Item(83, 9).Value = 6352.647
$ws.Cells.Item(83, 10).Value = 3481.818
$ws.Cells.Item(83, 11).Value = 31763.235
$ws.Cells.Item(83, 12).Value = 17409.09
$ws.Cells.Item(83, 13).Value = -26771.235
$ws.Cells.Item(83, 14).Value = -27393.09
$ws.Cells.Item(126, 8).Value = 2120
$ws.Cells.Item(126, 9).Value = 1600
$ws.Cells.Item(126, 10).Value = 3333.3333
$ws.Cells.Item(126, 11).Value = 4800
$ws.Cells.Item(126, 12).Value = 9999.999899999999
$ws.Cells.Item(126, 13).Value = -2330
$ws.Cells.Item(126, 14).Value = -14939.9999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 3222.8572
$ws.Cells.Item(82, 9).Value = 3533.3333
$ws.Cells.Item(82, 10).Value = 2990
$ws.Cells.Item(82, 11).Value = 3533.3333
$ws.Cells.Item(82, 12).Value = 2990
$ws.Cells.Item(82, 13).Value = -3172.3333
$ws.Cells.Item(82, 14).Value = -3712
$ws.Cells.Item(85, 8).Value = 3222.8572
$ws.Cells.Item(85, 9).Value = 3533.3333
$ws.Cells.Item(85, 10).Value = 2990
$ws.Cells.Item(85, 11).Value = 3533.3333
$ws.Cells.Item(85, 12).Value = 2990
$ws.Cells.Item(85, 13).Value = -2285.3333
$ws.Cells.Item(85, 14).Value = -5486
